$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old placeholder row 36 (B36 with leftover wrapText style / height)
$ws.Rows.Item(36).Delete()

# Row 27
$ws.Cells.Item(27, 1).Font.Size = 10
$ws.Cells.Item(27, 1).Value = 'poison_hemlock'
$ws.Cells.Item(27, 2).Font.Size = 10
$ws.Cells.Item(27, 2).Value = 'Characteristics'
$ws.Cells.Item(27, 3).Font.Size = 10
$ws.Cells.Item(27, 3).Value = 'Size'
$ws.Cells.Item(27, 4).Font.Size = 10
$ws.Cells.Item(27, 4).Value = '2 to 10 feet tall, with some sources stating it can reach up to 12 feet'

# Row 28
$ws.Cells.Item(28, 1).Font.Size = 10
$ws.Cells.Item(28, 1).Value = 'poison_hemlock'
$ws.Cells.Item(28, 2).Font.Size = 10
$ws.Cells.Item(28, 2).Value = 'Characteristics'
$ws.Cells.Item(28, 3).Font.Size = 10
$ws.Cells.Item(28, 3).Value = 'Shape'
$ws.Cells.Item(28, 4).Font.Size = 10
$ws.Cells.Item(28, 4).Value = 'triangular, lacy, fern-like leaf shape'

# Row 29
$ws.Cells.Item(29, 1).Font.Size = 10
$ws.Cells.Item(29, 1).Value = 'poison_hemlock'
$ws.Cells.Item(29, 2).Font.Size = 10
$ws.Cells.Item(29, 2).Value = 'Characteristics'
$ws.Cells.Item(29, 3).Font.Size = 10
$ws.Cells.Item(29, 3).Value = 'Color'
$ws.Cells.Item(29, 4).Font.Size = 10
$ws.Cells.Item(29, 4).Value = 'grayish-brown color when mature.'

# Row 30
$ws.Cells.Item(30, 1).Font.Size = 10
$ws.Cells.Item(30, 1).Value = 'poison_hemlock'
$ws.Cells.Item(30, 2).Font.Size = 10
$ws.Cells.Item(30, 2).Value = 'Characteristics'
$ws.Cells.Item(30, 3).Font.Size = 10
$ws.Cells.Item(30, 3).Value = 'body'
$ws.Cells.Item(30, 4).Font.Size = 10
$ws.Cells.Item(30, 4).Value = 'smooth, hollow, and ridged stems with distinctive purple spots, which are covered in hairless and lacy, fern-like leaves that resemble parsley'

# Row 31
$ws.Cells.Item(31, 1).Font.Size = 10
$ws.Cells.Item(31, 1).Value = 'poison_hemlock'
$ws.Cells.Item(31, 2).Font.Size = 10
$ws.Cells.Item(31, 2).Value = 'Health concern'
$ws.Cells.Item(31, 3).Font.Size = 10
$ws.Cells.Item(31, 3).Value = 'Risk factors'
$ws.Cells.Item(31, 4).Value = 'Central nervous system suppression, coma, kidney failure, low blood pressure, muscle breakdown, muscle death, muscle paralysis, adnd respiratory paralysis'

# Row 32
$ws.Cells.Item(32, 1).Font.Size = 10
$ws.Cells.Item(32, 1).Value = 'poison_hemlock'
$ws.Cells.Item(32, 2).Font.Size = 10
$ws.Cells.Item(32, 2).Value = 'What to do'
$ws.Cells.Item(32, 4).Font.Size = 10
$ws.Cells.Item(32, 4).Value = 'instantly seek emergency help, get fresh air, throw out all contaiminated things, wash your skin, and wash your eyes'

# Row 33
$ws.Cells.Item(33, 1).Font.Size = 10
$ws.Cells.Item(33, 1).Value = 'human_botfly'
$ws.Cells.Item(33, 2).Font.Size = 10
$ws.Cells.Item(33, 2).Value = 'Characteristics'
$ws.Cells.Item(33, 3).Font.Size = 10
$ws.Cells.Item(33, 3).Value = 'Size'
$ws.Cells.Item(33, 4).Value = '12-18mm long'

# Row 34
$ws.Cells.Item(34, 1).Font.Size = 10
$ws.Cells.Item(34, 1).Value = 'human_botfly'
$ws.Cells.Item(34, 2).Font.Size = 10
$ws.Cells.Item(34, 2).Value = 'Characteristics'
$ws.Cells.Item(34, 3).Font.Size = 10
$ws.Cells.Item(34, 3).Value = 'Shape'
$ws.Cells.Item(34, 4).Font.Size = 10
$ws.Cells.Item(34, 4).Value = 'A chunky, bee-like aperrence'

# Row 35
$ws.Cells.Item(35, 1).Font.Size = 10
$ws.Cells.Item(35, 1).Value = 'human_botfly'
$ws.Cells.Item(35, 2).Font.Size = 10
$ws.Cells.Item(35, 2).Value = 'Characteristics'
$ws.Cells.Item(35, 3).Font.Size = 10
$ws.Cells.Item(35, 3).Value = 'Color'
$ws.Cells.Item(35, 4).Font.Size = 10
$ws.Cells.Item(35, 4).Value = 'A yellow face, a metalic blue abdomen, and orange legs.'

# Row 36
$ws.Cells.Item(36, 1).Font.Size = 10
$ws.Cells.Item(36, 1).Value = 'human_botfly'
$ws.Cells.Item(36, 2).WrapText = $true
$ws.Cells.Item(36, 2).Value = 'Health concern'
$ws.Cells.Item(36, 3).Value = 'Pain and Discomfort'
$ws.Cells.Item(36, 4).Value = 'The developing larva creates a firm, painful, and itchy lump (known as a "warble") under the skin. A person may feel sensations of movement or sharp, stabbing pain as the larva moves or anchors itself with its spines.'

# Row 37
$ws.Cells.Item(37, 1).Font.Size = 10
$ws.Cells.Item(37, 1).Value = 'human_botfly'
$ws.Cells.Item(37, 2).WrapText = $true
$ws.Cells.Item(37, 2).Value = 'Health concern'
$ws.Cells.Item(37, 3).Value = 'Inflammation and Discharge'
$ws.Cells.Item(37, 4).Value = 'A local inflammatory response with redness, swelling, and a serous, bloody, or pus-like discharge from a central breathing hole (punctum) is common.'

# Row 38
$ws.Cells.Item(38, 1).Font.Size = 10
$ws.Cells.Item(38, 1).Value = 'human_botfly'
$ws.Cells.Item(38, 2).WrapText = $true
$ws.Cells.Item(38, 2).Value = 'Health concern'
$ws.Cells.Item(38, 3).Value = 'Tissue Damage and Scarring'
$ws.Cells.Item(38, 4).Value = 'As the larva grows, it consumes living tissue and can cause damage to the muscle. After the larva is removed and the wound heals, a significant scar may be left behind.'

# Row 39
$ws.Cells.Item(39, 1).Font.Size = 10
$ws.Cells.Item(39, 1).Value = 'German shepherd'
$ws.Cells.Item(39, 2).WrapText = $true
$ws.Cells.Item(39, 2).Value = 'Characteristics'
$ws.Cells.Item(39, 3).Value = 'Size'
$ws.Cells.Item(39, 4).Value = 'Males stand 24 to 26 inches tall at the shoulder and typically weigh between 65 to 90 pounds. Females are slightly smaller, standing 22 to 24 inches and weighing 50 to 70 pounds.'

# Row 40
$ws.Cells.Item(40, 1).Font.Size = 10
$ws.Cells.Item(40, 1).Value = 'German shepherd'
$ws.Cells.Item(40, 2).WrapText = $true
$ws.Cells.Item(40, 2).Value = 'Characteristics'
$ws.Cells.Item(40, 3).Value = 'Build'
$ws.Cells.Item(40, 4).Value = 'They have a strong, athletic build with a body that is typically longer than it is tall. The back can be straight (common in working lines) or sloped toward the hindquarters (often seen in show lines).'

# Row 41
$ws.Cells.Item(41, 1).Font.Size = 10
$ws.Cells.Item(41, 1).Value = 'German shepherd'
$ws.Cells.Item(41, 2).WrapText = $true
$ws.Cells.Item(41, 2).Value = 'Characteristics'
$ws.Cells.Item(41, 3).Value = 'Head and Face'
$ws.Cells.Item(41, 4).Value = 'The head is proportionate to the body, tapering to a strong, square-cut muzzle with a black nose and powerful jaws. They have medium-sized, dark, almond-shaped eyes that convey a keen and intelligent expression.'

# Row 42
$ws.Cells.Item(42, 1).Font.Size = 10
$ws.Cells.Item(42, 1).Value = 'German shepherd'
$ws.Cells.Item(42, 2).WrapText = $true
$ws.Cells.Item(42, 2).Value = 'Characteristics'
$ws.Cells.Item(42, 3).Value = 'Ears'
$ws.Cells.Item(42, 4).Value = 'A signature trait is their large, medium-sized ears which are broad at the base and carried erect and parallel when alert.'

# Row 43
$ws.Cells.Item(43, 1).Font.Size = 10
$ws.Cells.Item(43, 1).Value = 'German shepherd'
$ws.Cells.Item(43, 2).WrapText = $true
$ws.Cells.Item(43, 2).Value = 'Characteristics'
$ws.Cells.Item(43, 3).Value = 'Tail'
$ws.Cells.Item(43, 4).Value = 'They possess a long, bushy tail that hangs in a slight saber-like curve when at rest.'

# Row 44
$ws.Cells.Item(44, 1).Font.Size = 10
$ws.Cells.Item(44, 1).Value = 'German shepherd'
$ws.Cells.Item(44, 2).WrapText = $true
$ws.Cells.Item(44, 2).Value = 'Characteristics'
$ws.Cells.Item(44, 3).Value = 'Coat and Color'
$ws.Cells.Item(44, 4).Value = 'The German Shepherd has a dense double coat to protect them in various weather. The outer coat is typically medium-length, dense, and straight or slightly wavy. Common colors include black and tan/red, black and cream, sable, and solid black.'

# Row 45
$ws.Cells.Item(45, 1).Font.Size = 10
$ws.Cells.Item(45, 1).Value = 'German shepherd'
$ws.Cells.Item(45, 2).WrapText = $true
$ws.Cells.Item(45, 2).Value = 'Health concern'
$ws.Cells.Item(45, 3).Value = 'Deep Puncture Wounds and Lacerations'
$ws.Cells.Item(45, 4).Value = 'German Shepherds have strong jaws (bite force can exceed 238 pounds per square inch) capable of inflicting deep puncture wounds and significant soft tissue damage.'

# Row 46
$ws.Cells.Item(46, 1).Font.Size = 10
$ws.Cells.Item(46, 1).Value = 'German shepherd'
$ws.Cells.Item(46, 2).WrapText = $true
$ws.Cells.Item(46, 2).Value = 'Health concern'
$ws.Cells.Item(46, 3).Value = 'Crushing Injuries and Bone Fractures'
$ws.Cells.Item(46, 4).Value = 'The force of the bite can crush tissue and, in severe cases, break bones, especially in smaller victims like children or in sensitive areas like the hands and feet.'

# Row 47
$ws.Cells.Item(47, 1).Font.Size = 10
$ws.Cells.Item(47, 1).Value = 'German shepherd'
$ws.Cells.Item(47, 2).WrapText = $true
$ws.Cells.Item(47, 2).Value = 'Health concern'
$ws.Cells.Item(47, 3).Value = 'Nerve and Tendon Damage'
$ws.Cells.Item(47, 4).Value = 'The deep nature of the wounds can sever or damage nerves, tendons, and muscles, potentially leading to permanent loss of function, sensation, or chronic pain in the affected area.'

# Row 48
$ws.Cells.Item(48, 1).Font.Size = 10
$ws.Cells.Item(48, 1).Value = 'German shepherd'
$ws.Cells.Item(48, 2).WrapText = $true
$ws.Cells.Item(48, 2).Value = 'Health concern'
$ws.Cells.Item(48, 3).Value = 'Scarring and Disfigurement'
$ws.Cells.Item(48, 4).Value = 'Severe bites often require stitches or reconstructive surgery, which can result in permanent scarring or disfigurement.'

# Row 49
$ws.Cells.Item(49, 1).Font.Size = 10
$ws.Cells.Item(49, 1).Value = 'German shepherd'
$ws.Cells.Item(49, 2).WrapText = $true
$ws.Cells.Item(49, 2).Value = 'Health concern'
$ws.Cells.Item(49, 3).Value = 'Infectious Diseases and Complications'
$ws.Cells.Item(49, 4).Value = 'Dog mouths harbor numerous bacteria, and any bite that breaks the skin is at risk of infection.'

# Row 50
$ws.Cells.Item(50, 1).Font.Size = 10
$ws.Cells.Item(50, 1).Value = 'German shepherd'
$ws.Cells.Item(50, 2).WrapText = $true
$ws.Cells.Item(50, 2).Value = 'Health concern'
$ws.Cells.Item(50, 3).Value = 'Bacterial Infections'
$ws.Cells.Item(50, 4).Value = 'These are the most common health concerns.'

# Row 51
$ws.Cells.Item(51, 1).Font.Size = 10
$ws.Cells.Item(51, 1).Value = 'German shepherd'
$ws.Cells.Item(51, 2).WrapText = $true
$ws.Cells.Item(51, 2).Value = 'Health concern'
$ws.Cells.Item(51, 3).Value = 'Cellulitis'
$ws.Cells.Item(51, 4).Value = 'A rapidly spreading bacterial skin infection caused by bacteria likePasteurella multocida, which is common in dog saliva.'

# Row 52
$ws.Cells.Item(52, 1).Font.Size = 10
$ws.Cells.Item(52, 1).Value = 'German shepherd'
$ws.Cells.Item(52, 2).WrapText = $true
$ws.Cells.Item(52, 2).Value = 'Health concern'
$ws.Cells.Item(52, 3).Value = 'StaphandStrep Infections'
$ws.Cells.Item(52, 4).Value = 'Common bacterial infections that can lead to abscesses, joint swelling, and fever.'

# Row 53
$ws.Cells.Item(53, 1).Font.Size = 10
$ws.Cells.Item(53, 1).Value = 'German shepherd'
$ws.Cells.Item(53, 2).WrapText = $true
$ws.Cells.Item(53, 2).Value = 'Health concern'
$ws.Cells.Item(53, 3).Value = 'Capnocytophaga canimorsusinfection'
$ws.Cells.Item(53, 4).Value = 'A rare but potentially fatal infection, especially for individuals with weakened immune systems or who are asplenic (lacking a spleen). Symptoms can progress rapidly to sepsis, organ failure, andgangrene.'

# Row 54
$ws.Cells.Item(54, 1).Font.Size = 10
$ws.Cells.Item(54, 1).Value = 'German shepherd'
$ws.Cells.Item(54, 2).WrapText = $true
$ws.Cells.Item(54, 2).Value = 'Health concern'
$ws.Cells.Item(54, 3).Value = 'Sepsis'
$ws.Cells.Item(54, 4).Value = 'If an infection enters the bloodstream and spreads throughout the body, it can lead to sepsis, a life-threatening condition requiring immediate medical attention.'

# Row 55
$ws.Cells.Item(55, 1).Font.Size = 10
$ws.Cells.Item(55, 1).Value = 'German shepherd'
$ws.Cells.Item(55, 2).WrapText = $true
$ws.Cells.Item(55, 2).Value = 'Health concern'
$ws.Cells.Item(55, 3).Value = 'Rabies'
$ws.Cells.Item(55, 4).Value = 'This is a severe, almost always fatal viral disease once symptoms appear. The risk depends on the dog''s vaccination status and local rabies prevalence. If the dog is a stray or unvaccinated, post-exposure prophylaxis (PEP) is necessary.'

# Row 56
$ws.Cells.Item(56, 1).Font.Size = 10
$ws.Cells.Item(56, 1).Value = 'German shepherd'
$ws.Cells.Item(56, 2).WrapText = $true
$ws.Cells.Item(56, 2).Value = 'Health concern'
$ws.Cells.Item(56, 3).Value = 'Tetanus'
$ws.Cells.Item(56, 4).Value = 'The bacteria that cause tetanus can enter the body through a puncture wound. A tetanus shot booster may be needed if it has been more than five years since your last vaccination.'

# Row 57
$ws.Cells.Item(57, 1).Font.Size = 10
$ws.Cells.Item(57, 1).Value = 'German shepherd'
$ws.Cells.Item(57, 2).WrapText = $true
$ws.Cells.Item(57, 2).Value = 'Health concern'
$ws.Cells.Item(57, 3).Value = 'Psychological Effects:Post-Traumatic Stress Disorder (PTSD)'
$ws.Cells.Item(57, 4).Value = 'Victims may experience vivid flashbacks and intense fear.'

# Row 58
$ws.Cells.Item(58, 1).Font.Size = 10
$ws.Cells.Item(58, 1).Value = 'German shepherd'
$ws.Cells.Item(58, 2).WrapText = $true
$ws.Cells.Item(58, 2).Value = 'Health concern'
$ws.Cells.Item(58, 3).Value = 'Psychological Effects:Anxiety and Fear'
$ws.Cells.Item(58, 4).Value = 'Developing a lasting fear of dogs, which can impact daily life and social interactions.'

# Row 59
$ws.Cells.Item(59, 1).Font.Size = 10
$ws.Cells.Item(59, 1).Value = 'German shepherd'
$ws.Cells.Item(59, 2).WrapText = $true
$ws.Cells.Item(59, 2).Value = 'Health concern'
$ws.Cells.Item(59, 3).Value = 'Psychological Effects:Depression'
$ws.Cells.Item(59, 4).Value = 'The combination of physical injury and emotional distress can lead to depression.'

# Row 60
$ws.Cells.Item(60, 1).Font.Size = 10
$ws.Cells.Item(60, 1).Value = 'German shepherd'
$ws.Cells.Item(60, 2).Value = 'What to do'
$ws.Cells.Item(60, 4).Value = 'You should seek immediate medical care after any dog bite that breaks the skin. Prompt medical evaluation is critical to prevent severe complications, including'

# Row 61
$ws.Cells.Item(61, 1).Font.Size = 10
$ws.Cells.Item(61, 1).Value = 'German shepherd'
$ws.Cells.Item(61, 2).Value = 'What to do'
$ws.Cells.Item(61, 4).Value = 'Uncontrollable bleeding.'

# Row 62
$ws.Cells.Item(62, 1).Font.Size = 10
$ws.Cells.Item(62, 1).Value = 'German shepherd'
$ws.Cells.Item(62, 2).Value = 'What to do'
$ws.Cells.Item(62, 4).Value = 'A deep or large wound.'

# Row 63
$ws.Cells.Item(63, 1).Font.Size = 10
$ws.Cells.Item(63, 1).Value = 'German shepherd'
$ws.Cells.Item(63, 2).Value = 'What to do'
$ws.Cells.Item(63, 4).Value = 'Bites on the hands, face, neck, or feet.'

# Row 64
$ws.Cells.Item(64, 1).Font.Size = 10
$ws.Cells.Item(64, 1).Value = 'German shepherd'
$ws.Cells.Item(64, 2).Value = 'What to do'
$ws.Cells.Item(64, 4).Value = 'Signs of infection like spreading redness, swelling, increasing pain, or pus.'

# Row 65
$ws.Cells.Item(65, 1).Font.Size = 10
$ws.Cells.Item(65, 1).Value = 'German shepherd'
$ws.Cells.Item(65, 2).Value = 'What to do'
$ws.Cells.Item(65, 4).Value = 'If you are unsure of the dog''s vaccination status or have not had a tetanus shot in the last 5-10 years.'

# Update selection to match the saved view state
$ws.Range("D56").Select()